$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/19/2024  Through  8/25/2024"

# --- Helper functions for cells whose type (number <-> text) changes ---
function Set-NumFromDonor {
    param($ws, $addr, $donor, $value)
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($addr).Value = $value
}
function Set-TextFromDonor {
    param($ws, $addr, $donor, $text)
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("ZZ1").NumberFormat = "@"
    $ws.Range("ZZ1").Value = $text
    $ws.Range("ZZ1").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# --- Cells whose type/style changes (number <-> text N/A placeholder) ---
Set-NumFromDonor $ws "D14" "D16" 1
Set-NumFromDonor $ws "E14" "E16" -100
Set-TextFromDonor $ws "C15" "C29" "0"
Set-NumFromDonor $ws "D18" "D16" 1
Set-NumFromDonor $ws "E18" "E16" 600
Set-TextFromDonor $ws "C22" "C29" "0"
Set-NumFromDonor $ws "C23" "C16" 3
Set-TextFromDonor $ws "D23" "D29" "0"
Set-TextFromDonor $ws "E23" "E29" "***.*"
Set-NumFromDonor $ws "C28" "C16" 1
Set-NumFromDonor $ws "C33" "C16" 2
Set-NumFromDonor $ws "F33" "F16" 2
Set-NumFromDonor $ws "I33" "I16" 2

$ws.Range("ZZ1").Clear() | Out-Null

# --- Plain value updates (style unchanged) ---
$ws.Range("G14").Value = 2
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = 0
$ws.Range("L15").Value = -45.454545454545
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = -54.545454545454
$ws.Range("I16").Value = 107
$ws.Range("J16").Value = 109
$ws.Range("K16").Value = -1.834862385321
$ws.Range("L16").Value = -17.692307692307
$ws.Range("M16").Value = -29.605263157894
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 162
$ws.Range("J17").Value = 164
$ws.Range("K17").Value = -1.219512195121
$ws.Range("L17").Value = -14.736842105263
$ws.Range("M17").Value = 32.786885245901
$ws.Range("C18").Value = 7
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 125
$ws.Range("I18").Value = 71
$ws.Range("J18").Value = 77
$ws.Range("K18").Value = -7.792207792207
$ws.Range("L18").Value = 10.9375
$ws.Range("M18").Value = 9.230769230769
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -45.454545454545
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = -38.095238095238
$ws.Range("I19").Value = 228
$ws.Range("J19").Value = 235
$ws.Range("K19").Value = -2.978723404255
$ws.Range("L19").Value = 13.432835820895
$ws.Range("M19").Value = 32.558139534883
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 71
$ws.Range("J20").Value = 99
$ws.Range("K20").Value = -28.282828282828
$ws.Range("L20").Value = -24.468085106383
$ws.Range("M20").Value = 57.777777777777
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -4.347826086956
$ws.Range("F21").Value = 69
$ws.Range("G21").Value = 104
$ws.Range("H21").Value = -33.653846153846
$ws.Range("I21").Value = 650
$ws.Range("J21").Value = 695
$ws.Range("K21").Value = -6.474820143884
$ws.Range("L21").Value = -6.069364161849
$ws.Range("M21").Value = 14.235500878734
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -80
$ws.Range("J22").Value = 12
$ws.Range("K22").Value = -33.333333333333
$ws.Range("L22").Value = -61.904761904761
$ws.Range("M22").Value = -11.111111111111
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 50
$ws.Range("I23").Value = 21
$ws.Range("K23").Value = 133.333333333333
$ws.Range("L23").Value = 40
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -50
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = -1.149425287356
$ws.Range("I24").Value = 591
$ws.Range("J24").Value = 633
$ws.Range("K24").Value = -6.635071090047
$ws.Range("L24").Value = -37.394067796610
$ws.Range("M24").Value = 60.162601626016
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -71.428571428571
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 90
$ws.Range("I25").Value = 189
$ws.Range("J25").Value = 157
$ws.Range("K25").Value = 20.382165605095
$ws.Range("L25").Value = -66.900175131348
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 150
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = -43.243243243243
$ws.Range("I26").Value = 221
$ws.Range("J26").Value = 262
$ws.Range("K26").Value = -15.648854961832
$ws.Range("L26").Value = -19.343065693430
$ws.Range("M26").Value = -31.366459627329
$ws.Range("F27").Value = 2
$ws.Range("I27").Value = 10
$ws.Range("K27").Value = 11.111111111111
$ws.Range("L27").Value = -41.176470588235
$ws.Range("E28").Value = 0
$ws.Range("I28").Value = 29
$ws.Range("J28").Value = 26
$ws.Range("K28").Value = 11.538461538461
$ws.Range("L28").Value = -19.444444444444
$ws.Range("L29").Value = -66.666666666666
$ws.Range("L30").Value = -80
$ws.Range("H33").Value = 100
$ws.Range("K33").Value = 0
